$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ44951355",
    "summ46124965",
    "summ46978517",
    "summ47333171",
    "summ47671051",
    "summ48020392",
    "summ48346527",
    "summ48695839",
    "summ49070491"
)

for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
